$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 78.94643366666666
$ws.Range("H2").Value = 236.839301
$ws.Range("I2").Value = 0.8231230137654455
$ws.Range("J2").Value = 0.8231230137654454
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.807091
$ws.Range("N2").Value = 23.421273
$ws.Range("O2").Value = 0.3133326987643095
$ws.Range("P2").Value = 0.3133326987643095
$ws.Range("Q2").Value = 616.3419917611303
$ws.Range("R2").Value = 5547.077925850173
$ws.Range("S2").Value = 0.2579113553181389
$ws.Range("T2").Value = 0.2579113553181389

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 78.94643366666666
$ws.Range("H3").Value = 236.839301
$ws.Range("I3").Value = 0.8231230137654455
$ws.Range("J3").Value = 0.8231230137654454
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 12.70064766666667
$ws.Range("N3").Value = 38.101943
$ws.Range("O3").Value = 0.5097325251430138
$ws.Range("P3").Value = 0.5097325251430138
$ws.Range("Q3").Value = 1002.670838540205
$ws.Range("R3").Value = 9024.037546861842
$ws.Range("S3").Value = 0.4195725723099882
$ws.Range("T3").Value = 0.4195725723099882

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 78.94643366666666
$ws.Range("H4").Value = 236.839301
$ws.Range("I4").Value = 0.8231230137654455
$ws.Range("J4").Value = 0.8231230137654454
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.408559666666666
$ws.Range("N4").Value = 13.225679
$ws.Range("O4").Value = 0.1769347760926767
$ws.Range("P4").Value = 0.1769347760926767
$ws.Range("Q4").Value = 348.040063290042
$ws.Range("R4").Value = 3132.360569610379
$ws.Range("S4").Value = 0.1456390861373183
$ws.Range("T4").Value = 0.1456390861373183

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.569858333333333
$ws.Range("H5").Value = 10.709575
$ws.Range("I5").Value = 0.03722058633396773
$ws.Range("J5").Value = 0.03722058633396773
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.807091
$ws.Range("N5").Value = 23.421273
$ws.Range("O5").Value = 0.3133326987643095
$ws.Range("P5").Value = 0.3133326987643095
$ws.Range("Q5").Value = 27.87020886544166
$ws.Range("R5").Value = 250.831879788975
$ws.Range("S5").Value = 0.01166242676561208
$ws.Range("T5").Value = 0.01166242676561208

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.569858333333333
$ws.Range("H6").Value = 10.709575
$ws.Range("I6").Value = 0.03722058633396773
$ws.Range("J6").Value = 0.03722058633396773
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 12.70064766666667
$ws.Range("N6").Value = 38.101943
$ws.Range("O6").Value = 0.5097325251430138
$ws.Range("P6").Value = 0.5097325251430138
$ws.Range("Q6").Value = 45.33951291158055
$ws.Range("R6").Value = 408.055616204225
$ws.Range("S6").Value = 0.01897254345931692
$ws.Range("T6").Value = 0.01897254345931692

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.569858333333333
$ws.Range("H7").Value = 10.709575
$ws.Range("I7").Value = 0.03722058633396773
$ws.Range("J7").Value = 0.03722058633396773
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.408559666666666
$ws.Range("N7").Value = 13.225679
$ws.Range("O7").Value = 0.1769347760926767
$ws.Range("P7").Value = 0.1769347760926767
$ws.Range("Q7").Value = 15.73793346404722
$ws.Range("R7").Value = 141.641401176425
$ws.Range("S7").Value = 0.006585616109038721
$ws.Range("T7").Value = 0.006585616109038721

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 13.394565
$ws.Range("H8").Value = 40.183695
$ws.Range("I8").Value = 0.1396563999005869
$ws.Range("J8").Value = 0.1396563999005868
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.807091
$ws.Range("N8").Value = 23.421273
$ws.Range("O8").Value = 0.3133326987643095
$ws.Range("P8").Value = 0.3133326987643095
$ws.Range("Q8").Value = 104.572587860415
$ws.Range("R8").Value = 941.1532907437349
$ws.Range("S8").Value = 0.04375891668055852
$ws.Range("T8").Value = 0.04375891668055851

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 13.394565
$ws.Range("H9").Value = 40.183695
$ws.Range("I9").Value = 0.1396563999005869
$ws.Range("J9").Value = 0.1396563999005868
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 12.70064766666667
$ws.Range("N9").Value = 38.101943
$ws.Range("O9").Value = 0.5097325251430138
$ws.Range("P9").Value = 0.5097325251430138
$ws.Range("Q9").Value = 170.119650713265
$ws.Range("R9").Value = 1531.076856419385
$ws.Range("S9").Value = 0.07118740937370868
$ws.Range("T9").Value = 0.07118740937370867

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 13.394565
$ws.Range("H10").Value = 40.183695
$ws.Range("I10").Value = 0.1396563999005869
$ws.Range("J10").Value = 0.1396563999005868
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.408559666666666
$ws.Range("N10").Value = 13.225679
$ws.Range("O10").Value = 0.1769347760926767
$ws.Range("P10").Value = 0.1769347760926767
$ws.Range("Q10").Value = 59.05073901154499
$ws.Range("R10").Value = 531.456651103905
$ws.Range("S10").Value = 0.02471007384631965
$ws.Range("T10").Value = 0.02471007384631965
